# The workbook originally has three sheets, in order:
#   "UserList"      (sheetId 1) - 4 users (standard, locked_out, problem, performance_glitch)
#   "Sorting"       (sheetId 2) - sort-by test data
#   "UserList (2)"  (sheetId 4) - 3 users (standard, problem, performance_glitch)
#
# The edit renames/reorders the sheets so that:
#   - the old "UserList (2)" sheet becomes the first sheet, renamed to "UserList"
#   - the old "UserList" sheet becomes the second sheet, renamed to "UserList-1"
#   - the "Sorting" sheet becomes the last (third) sheet
# and makes the (new, first) "UserList" sheet the active sheet, with cell F6 selected.

$wb = $excel.ActiveWorkbook

$oldUserList    = $wb.Worksheets.Item("UserList")
$oldUserList2   = $wb.Worksheets.Item("UserList (2)")

# Rename the original "UserList" out of the way first, then promote
# "UserList (2)" to the now-free "UserList" name.
$oldUserList.Name = "UserList-1"
$oldUserList2.Name = "UserList"

# Move the (renamed) former "UserList (2)" sheet to the front of the workbook.
$oldUserList2.Move($oldUserList)

# Select the new "UserList" sheet and put the selection on cell F6, which also
# makes it the active/visible tab.
$newFirst = $wb.Worksheets.Item("UserList")
$newFirst.Activate() | Out-Null
$newFirst.Range("F6").Select() | Out-Null
